$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 116
$ws.Cells.Item(116, 1).Value = 8
$ws.Cells.Item(116, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(116, 3).Value = "Coquimbo"
$ws.Cells.Item(116, 4).Value = 44939
$ws.Cells.Item(116, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(116, 5).Value = 4
$ws.Cells.Item(116, 6).Value = "Fruta"
$ws.Cells.Item(116, 7).Value = 100103
$ws.Cells.Item(116, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(116, 9).Value = 100103003
$ws.Cells.Item(116, 10).Value = "Damasco"
$ws.Cells.Item(116, 11).Value = "Dina"
$ws.Cells.Item(116, 12).Value = "Especial"
$ws.Cells.Item(116, 13).Value = 200
$ws.Cells.Item(116, 14).Value = 22000
$ws.Cells.Item(116, 15).Value = 23000
$ws.Cells.Item(116, 16).Value = 22500
$ws.Cells.Item(116, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(116, 18).Value = "Región Metropolitana"
$ws.Cells.Item(116, 19).Value = 1406
$ws.Cells.Item(116, 20).Value = 16

# Row 117
$ws.Cells.Item(117, 1).Value = 8
$ws.Cells.Item(117, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(117, 3).Value = "Coquimbo"
$ws.Cells.Item(117, 4).Value = 44939
$ws.Cells.Item(117, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(117, 5).Value = 4
$ws.Cells.Item(117, 6).Value = "Fruta"
$ws.Cells.Item(117, 7).Value = 100103
$ws.Cells.Item(117, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(117, 9).Value = 100103003
$ws.Cells.Item(117, 10).Value = "Damasco"
$ws.Cells.Item(117, 11).Value = "Dina"
$ws.Cells.Item(117, 12).Value = "Primera"
$ws.Cells.Item(117, 13).Value = 240
$ws.Cells.Item(117, 14).Value = 20000
$ws.Cells.Item(117, 15).Value = 21000
$ws.Cells.Item(117, 16).Value = 20500
$ws.Cells.Item(117, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(117, 18).Value = "Región Metropolitana"
$ws.Cells.Item(117, 19).Value = 1281
$ws.Cells.Item(117, 20).Value = 16
